# Auto-generated edit script applying the Typhon_Profits market-data refresh diff
# Updates currentAveragePrice* / LevePrice* / LeveProfit* columns (H:N) across sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 2558
$ws.Range("I19").Value = 6883.6665
$ws.Range("J19").Value = 395.16666
$ws.Range("K19").Value = 6883.6665
$ws.Range("L19").Value = 395.16666
$ws.Range("M19").Value = -6708.6665
$ws.Range("N19").Value = -745.16666
$ws.Range("H33").Value = 214.3
$ws.Range("I33").Value = 214.3
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 214.3
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = 14.69999999999999
$ws.Range("N33").Value = ""
$ws.Range("H58").Value = 3735.2856
$ws.Range("I58").Value = 376.66666
$ws.Range("J58").Value = 6254.25
$ws.Range("K58").Value = 1129.99998
$ws.Range("L58").Value = 18762.75
$ws.Range("M58").Value = -979.99998
$ws.Range("N58").Value = -19062.75
$ws.Range("H82").Value = 1064
$ws.Range("I82").Value = 1064
$ws.Range("K82").Value = 3192
$ws.Range("M82").Value = -2786
$ws.Range("H85").Value = 1064
$ws.Range("I85").Value = 1064
$ws.Range("K85").Value = 3192
$ws.Range("M85").Value = -1788
$ws.Range("H98").Value = 615.8947
$ws.Range("I98").Value = 622.3333
$ws.Range("J98").Value = 500
$ws.Range("K98").Value = 622.3333
$ws.Range("L98").Value = 500
$ws.Range("M98").Value = 875.6667
$ws.Range("N98").Value = -3496
$ws.Range("H113").Value = 35717800
$ws.Range("I113").Value = 83336720
$ws.Range("J113").Value = 3615.8125
$ws.Range("K113").Value = 83336720
$ws.Range("L113").Value = 3615.8125
$ws.Range("M113").Value = -83333466
$ws.Range("N113").Value = -10123.8125
$ws.Range("H122").Value = 615.8947
$ws.Range("I122").Value = 622.3333
$ws.Range("J122").Value = 500
$ws.Range("K122").Value = 1866.9999
$ws.Range("L122").Value = 1500
$ws.Range("M122").Value = 583.0001
$ws.Range("N122").Value = -6400
$ws.Range("H129").Value = 625634.6
$ws.Range("J129").Value = 909896.56
$ws.Range("L129").Value = 2729689.68
$ws.Range("N129").Value = -2739689.68
$ws.Range("H138").Value = 2402.0286
$ws.Range("I138").Value = 3399.625
$ws.Range("K138").Value = 10198.875
$ws.Range("M138").Value = -5058.875
$ws.Range("H141").Value = 2280.8333
$ws.Range("I141").Value = 1807.2222
$ws.Range("J141").Value = 3701.6667
$ws.Range("K141").Value = 5421.6666
$ws.Range("L141").Value = 11105.0001
$ws.Range("M141").Value = -241.6665999999996
$ws.Range("N141").Value = -21465.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6314.6123
$ws.Range("I32").Value = 4661.2173
$ws.Range("J32").Value = 31666.666
$ws.Range("K32").Value = 4661.2173
$ws.Range("L32").Value = 31666.666
$ws.Range("M32").Value = -4374.2173
$ws.Range("N32").Value = -32240.666
$ws.Range("H74").Value = 43479090
$ws.Range("I74").Value = 71429176
$ws.Range("K74").Value = 71429176
$ws.Range("M74").Value = -71428302
$ws.Range("H77").Value = 43479090
$ws.Range("I77").Value = 71429176
$ws.Range("K77").Value = 357145880
$ws.Range("M77").Value = -357141512
$ws.Range("H110").Value = 829.36365
$ws.Range("I110").Value = 712.3
$ws.Range("K110").Value = 712.3
$ws.Range("M110").Value = 1332.7
$ws.Range("H122").Value = 3532.111
$ws.Range("I122").Value = 3348.625
$ws.Range("K122").Value = 10045.875
$ws.Range("M122").Value = -7595.875
$ws.Range("H132").Value = 18739.934
$ws.Range("I132").Value = 1986.3636
$ws.Range("J132").Value = 64812.25
$ws.Range("K132").Value = 5959.0908
$ws.Range("L132").Value = 194436.75
$ws.Range("M132").Value = -3429.0908
$ws.Range("N132").Value = -199496.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H60").Value = 0
$ws.Range("J60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("N60").Value = ""
$ws.Range("H82").Value = 20012.1
$ws.Range("I82").Value = 6669.4287
$ws.Range("J82").Value = 51145
$ws.Range("K82").Value = 6669.4287
$ws.Range("L82").Value = 51145
$ws.Range("M82").Value = -6286.4287
$ws.Range("N82").Value = -51911
$ws.Range("H85").Value = 20012.1
$ws.Range("I85").Value = 6669.4287
$ws.Range("J85").Value = 51145
$ws.Range("K85").Value = 6669.4287
$ws.Range("L85").Value = 51145
$ws.Range("M85").Value = -5343.4287
$ws.Range("N85").Value = -53797
$ws.Range("H110").Value = 45700
$ws.Range("J110").Value = 45700
$ws.Range("L110").Value = 45700
$ws.Range("N110").Value = -53880
$ws.Range("H134").Value = 4269.4
$ws.Range("I134").Value = 4647.0454
$ws.Range("J134").Value = 1500
$ws.Range("K134").Value = 13941.1362
$ws.Range("L134").Value = 4500
$ws.Range("M134").Value = -11406.1362
$ws.Range("N134").Value = -9570

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1401.875
$ws.Range("I122").Value = 1950
$ws.Range("J122").Value = 1219.1666
$ws.Range("K122").Value = 5850
$ws.Range("L122").Value = 3657.4998
$ws.Range("M122").Value = -3400
$ws.Range("N122").Value = -8557.4998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 204.4
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 204.4
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 613.2
$ws.Range("M12").Value = ""
$ws.Range("N12").Value = -959.2
$ws.Range("H75").Value = 2225.5
$ws.Range("I75").Value = 1847
$ws.Range("J75").Value = 2604
$ws.Range("K75").Value = 5541
$ws.Range("L75").Value = 7812
$ws.Range("M75").Value = -4543
$ws.Range("N75").Value = -9808
$ws.Range("H78").Value = 2225.5
$ws.Range("I78").Value = 1847
$ws.Range("J78").Value = 2604
$ws.Range("K78").Value = 16623
$ws.Range("L78").Value = 23436
$ws.Range("M78").Value = -11631
$ws.Range("N78").Value = -33420
$ws.Range("H103").Value = 491.81818
$ws.Range("I103").Value = 454.8889
$ws.Range("J103").Value = 658
$ws.Range("K103").Value = 1364.6667
$ws.Range("L103").Value = 1974
$ws.Range("M103").Value = -485.6667
$ws.Range("N103").Value = -3732
$ws.Range("H131").Value = 747.59
$ws.Range("J131").Value = 747.61615
$ws.Range("L131").Value = 2242.84845
$ws.Range("N131").Value = -12322.84845
$ws.Range("H132").Value = 1174.75
$ws.Range("I132").Value = 470.5
$ws.Range("J132").Value = 1409.5
$ws.Range("K132").Value = 4234.5
$ws.Range("L132").Value = 12685.5
$ws.Range("M132").Value = -1704.5
$ws.Range("N132").Value = -17745.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1449.75
$ws.Range("I122").Value = 999.6667
$ws.Range("J122").Value = 2800
$ws.Range("K122").Value = 2999.0001
$ws.Range("L122").Value = 8400
$ws.Range("M122").Value = -549.0001000000002
$ws.Range("N122").Value = -13300
$ws.Range("H132").Value = 26196.043
$ws.Range("I132").Value = 5068.8125
$ws.Range("J132").Value = 74486.86
$ws.Range("K132").Value = 15206.4375
$ws.Range("L132").Value = 223460.58
$ws.Range("M132").Value = -12676.4375
$ws.Range("N132").Value = -228520.58

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 853.875
$ws.Range("I46").Value = 863.6429
$ws.Range("J46").Value = 840.2
$ws.Range("K46").Value = 863.6429
$ws.Range("L46").Value = 840.2
$ws.Range("M46").Value = -675.6429
$ws.Range("N46").Value = -1216.2
$ws.Range("H55").Value = 72.21429
$ws.Range("I55").Value = 67
$ws.Range("K55").Value = 67
$ws.Range("M55").Value = 106
$ws.Range("H82").Value = 5334.3335
$ws.Range("I82").Value = 4500
$ws.Range("J82").Value = 7003
$ws.Range("K82").Value = 4500
$ws.Range("L82").Value = 7003
$ws.Range("M82").Value = -4139
$ws.Range("N82").Value = -7725
$ws.Range("H85").Value = 5334.3335
$ws.Range("I85").Value = 4500
$ws.Range("J85").Value = 7003
$ws.Range("K85").Value = 4500
$ws.Range("L85").Value = 7003
$ws.Range("M85").Value = -3252
$ws.Range("N85").Value = -9499
$ws.Range("H122").Value = 1092197.4
$ws.Range("I122").Value = 1963384.8
$ws.Range("J122").Value = 3213.125
$ws.Range("K122").Value = 5890154.4
$ws.Range("L122").Value = 9639.375
$ws.Range("M122").Value = -5887704.4
$ws.Range("N122").Value = -14539.375
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").Value = ""
$ws.Range("H132").Value = 2720.7144
$ws.Range("I132").Value = 1948.8334
$ws.Range("J132").Value = 3299.625
$ws.Range("K132").Value = 5846.5002
$ws.Range("L132").Value = 9898.875
$ws.Range("M132").Value = -3316.5002
$ws.Range("N132").Value = -14958.875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 3374.75
$ws.Range("J15").Value = 3374.75
$ws.Range("L15").Value = 3374.75
$ws.Range("N15").Value = -3950.75
$ws.Range("H132").Value = 917.8222
$ws.Range("I132").Value = 655.6857
$ws.Range("J132").Value = 1835.3
$ws.Range("K132").Value = 1967.0571
$ws.Range("L132").Value = 5505.9
$ws.Range("M132").Value = 562.9429
$ws.Range("N132").Value = -10565.9

